$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("M")

# Insert a new row above current row 10 (shifts existing rows down)
$ws.Rows.Item(10).Insert()

$cell = $ws.Range("A10")
$cell.Value = "Groundwater Management (SWRCB)"
$cell.Font.Name = "Calibri"

$ws.Activate()
$ws.Range("D15").Select()
